$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.216.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.025.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.81"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.10"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0782"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.325.61"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.29"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.12"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.738"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.17"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.025.15"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.167.62"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.16"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.07"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.42"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.66"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.05"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.57%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.72"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.31"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.118"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.43%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.62%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.20"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "BinanceUSD"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.49"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.474.81"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0214"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "94.59"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0914"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.25"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +13.92%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.82%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.16"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.212.32"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.04%  "
